# EquipmentDisposition.docx update
#
# The "Online Photography Retailers" section's B&H entry ends with a
# sentence that names a URL as plain text:
#   "...Their hours of operation can be found at:
#    https://www.bhphotovideo.com/find/HelpCenter/StoreInfo.jsp."
#
# Turn that plain-text URL into a real hyperlink (same treatment the
# other retailer/museum links in the document already get), leaving the
# rest of the sentence - and the trailing period after the URL - as
# plain text.

$d = $word.ActiveDocument

$url = "https://www.bhphotovideo.com/find/HelpCenter/StoreInfo.jsp"

# Find the URL text (the trailing "." belongs to the sentence, not the
# address, so it is deliberately excluded from the search string).
$urlRange = $d.Content
$found = $urlRange.Find.Execute($url, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "", 0)

if ($found) {
    # Converts the plain-text range into a hyperlink pointing at the
    # same address, displaying the original URL text.
    $d.Hyperlinks.Add($urlRange, $url) | Out-Null
} else {
    Write-Host "WARNING: could not locate the B&H hours-of-operation URL"
}
